$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("C9:C10")
$texts = @("DAC","ADC","D1V5","D3V3","VSUP","DGND","A3V3","A1V5","AGND")
foreach ($t in $texts) {
    $fc = $rng.FormatConditions.Add(9, 0, $null, $null, $t)
    $fc.Interior.Color() = 13561798
}
Write-Host "done"
